$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 10 ("Weakness of Our Project" -> "Risk of Our Project")
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)

# Title text change
$s10.Shapes.Item(1).TextFrame.TextRange.Text = "Risk of Our Project"

# Body placeholder: replace the single placeholder space with four bullet
# paragraphs (they all inherit the existing Wingdings "q" bullet pPr that is
# already on the first paragraph).
$body10 = $s10.Shapes.Item(2).TextFrame.TextRange
$body10.Text = "Node.js is recent technology. So outcoming problem handle is more challenge.`rThere is no good community scope.`rThe information of house advertiser keep secure is more challenge.`rSometimes bad people to take information, subscription is apply on house advertiser."

# ---------------------------------------------------------------------------
# Slide 5 ("Feasibility Study") - fix the "Environmenttal" typo
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$body5 = $s5.Shapes.Item(2).TextFrame.TextRange
$full5 = $body5.Text
$idx5 = $full5.IndexOf("Environmenttal")
if ($idx5 -ge 0) {
    $sub5 = $body5.Characters($idx5 + 1, "Environmenttal".Length)
    $sub5.Text = "Environmental"
}

# ---------------------------------------------------------------------------
# Slide 9 ("Expected Outcome") - add the "Complete business Platform." bullet
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$body9 = $s9.Shapes.Item(2).TextFrame.TextRange
$body9.Text = "Complete business Platform.`r"

# Turn the plain ("no bullet") paragraph format into the Wingdings "q" bullet
# used throughout the rest of the deck.
$pf9 = $body9.ParagraphFormat
$bullet9 = $pf9.Bullet
$bullet9.Visible = -1
$bullet9.Character = 113
$bullet9.Font.Name = "Wingdings"
